# Bugfixed evaluation and simulated rt_data for components: refresh the forecast
# vector rows (date, y_0, y_0_forecast, y_1, y_1_forecast) and append the newest row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 is new. Copy the date-number-format/border/bold style from A52 into A53
# first so the new date cell keeps the same "s=2" style as the rest of column A.
$ws.Cells.Item(52, 1).Copy($ws.Cells.Item(53, 1))

# Columns: A=date_of_forecast, B=y_0, C=y_0_forecast, D=y_1, E=y_1_forecast
$rows = @(
  @(2, 39400, 2007, 5.896808312953783, 2008, 10.70171490310616),
  @(3, 39583, 2008, 2.316462327490609, 2009, -12.90434384910116),
  @(4, 39765, 2008, 7.441962824572235, 2009, 10.14224506046018),
  @(5, 39948, 2009, 9.162319139560072, 2010, 11.00161740362515),
  @(6, 40130, 2009, 6.277541464866987, 2010, 9.926356894615008),
  @(7, 40310, 2010, 4.440675708811326, 2011, 0.02592873473603596),
  @(8, 40494, 2010, 6.535114773304773, 2011, 6.670100100023779),
  @(9, 40676, 2011, 6.166830030198267, 2012, 5.009734572843771),
  @(10, 40862, 2011, 5.12051970717502, 2012, 6.878526556821973),
  @(11, 41044, 2012, 2.204591531842581, 2013, -4.726949348040732),
  @(12, 41228, 2012, 3.65682115264816, 2013, 4.783388268072009),
  @(13, 41409, 2013, 1.09689030557385, 2014, -4.442195584720931),
  @(14, 41592, 2013, 2.943878639034381, 2014, 5.440628060653574),
  @(15, 41774, 2014, 1.600060471414833, 2015, -2.949819494134909),
  @(16, 41957, 2014, 1.172679597477866, 2015, 0.7846309372076199),
  @(17, 42137, 2015, 1.651055586686678, 2016, 2.833734372666652),
  @(18, 42321, 2015, 2.961845079861303, 2016, 2.521390070591267),
  @(19, 42503, 2016, 1.99440460461342, 2017, 0.9898525518331924),
  @(20, 42689, 2016, 2.508469427909898, 2017, 2.540874511056623),
  @(21, 42867, 2017, 2.613530175870626, 2018, 3.313641510414356),
  @(22, 43053, 2017, 3.523703831572056, 2018, 3.025650759929999),
  @(23, 43145, 2018, 1.274170701304334, 2019, -0.2174071002483902),
  @(24, 43235, 2018, 3.696331036365752, 2019, 6.233648892987009),
  @(25, 43326, 2018, 1.207340632843712, 2019, 1.091965090482971),
  @(26, 43418, 2018, 1.178605266817589, 2019, 0.6761281928490348),
  @(27, 43510, 2019, 1.37810805366898, 2020, 1.937142119249913),
  @(28, 43600, 2019, -0.973865931199458, 2020, -4.327930935899992),
  @(29, 43691, 2019, 3.031664248068777, 2020, 3.470008002361924),
  @(30, 43783, 2019, 3.047037961814492, 2020, 3.638425237116749),
  @(31, 43875, 2020, 3.770384266127369, 2021, 4.270473305078393),
  @(32, 43966, 2020, 1.725130460355095, 2021, -1.194610791900008),
  @(33, 44068, 2020, -0.4891791466461126, 2021, 1.831864451486243),
  @(34, 44159, 2020, -0.2228847697281378, 2021, 2.646620531497335),
  @(35, 44251, 2021, 2.288807135633641, 2022, 2.832246610071509),
  @(36, 44341, 2021, -0.5115004854862049, 2022, 2.23509962177757),
  @(37, 44432, 2021, -0.7941560676977599, 2022, 4.94197512951382),
  @(38, 44525, 2021, -1.165854108406617, 2022, 4.306257314632855),
  @(39, 44617, 2022, 4.061815910726119, 2023, 5.709386381624437),
  @(40, 44706, 2022, 0.9727820482463123, 2023, 2.866869504079239),
  @(41, 44798, 2022, 2.34069710769782, 2023, 1.687287030004092),
  @(42, 44890, 2022, 2.501311189006916, 2023, 3.687525625406263),
  @(43, 44981, 2023, 2.830392552592231, 2024, 2.56445306171762),
  @(44, 45071, 2023, 2.122313752051319, 2024, 0.8382457967197388),
  @(45, 45163, 2023, 0.9259311313598806, 2024, 0.9881529557394186),
  @(46, 45254, 2023, 0.6753076481029074, 2024, 0.5720316833079497),
  @(47, 45345, 2024, 1.488849104696177, 2025, 2.559682772914296),
  @(48, 45436, 2024, -0.1133034947815914, 2025, 0.8323378752418176),
  @(49, 45534, 2024, 1.935025917091848, 2025, 1.300166298184058),
  @(50, 45618, 2024, 2.039329803030121, 2025, 2.838786065810939),
  @(51, 45713, 2025, 1.774603592091717, 2026, 2.333011795730777),
  @(52, 45800, 2025, 3.161594928268019, 2026, 2.57979941834241),
  @(53, 45891, 2025, 2.43119486791763, 2026, 3.197192111131564)
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
}
